$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append to the results sheet (rows 229-242)
$data = @(
    @(39.78, 0.16, 43.16329216957092, 40.296875),
    @(37.37, 0.36, 55.52490019798279, 51.75),
    @(39.78, 0.16, 38.23706722259521, 34.296875),
    @(37.37, 0.36, 80.1137330532074, 74.453125),
    @(39.78, 0.16, 23.20184087753296, 20.78125),
    @(66.45999999999999, 0.6400000000000001, 168.7978188991547, 162.71875),
    @(69.63, 0.16, 43.12589979171753, 42.796875),
    @(37.37, 0.36, 82.05421042442322, 80.046875),
    @(39.78, 0.16, 34.1039719581604, 34.140625),
    @(39.78, 0.16, 38.44294476509094, 38.5625),
    @(39.78, 0.16, 38.57809257507324, 38.046875),
    @(39.78, 0.16, 20.44303488731384, 20.359375),
    @(37.37, 0.36, 45.1411452293396, 44.828125),
    @(39.78, 0.16, 28.77275109291077, 28.703125)
)

$startRow = 229
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
